$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Words for rows 10-19 (these become shared-string indices 0-9 because
# they are written first, and the engine appends new shared strings in
# the order the cell values are assigned).
$group1 = @("Plethora","Profuse","Prolific","Rife","Spate","Steeped","Surfeit","Surge","Teeming","Volley")
for ($i = 0; $i -lt $group1.Length; $i++) {
    $row = $i + 10
    $ws.Cells.Item($row, 1).Value = $group1[$i]
}

# Words for rows 1-9 (these become shared-string indices 10-18).
$group2 = @("Exiguous","Marginal","Meagre","Negligible","Paltry","Scanty","Skimpy","Spare","Sparse")
for ($i = 0; $i -lt $group2.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $group2[$i]
}

# Rows 20-22 no longer have a word in column A.
$ws.Cells.Item(20, 1).Value = ""
$ws.Cells.Item(21, 1).Value = ""
$ws.Cells.Item(22, 1).Value = ""

# Update the selection shown in the sheet view to D1.
$ws.Range("D1").Select() | Out-Null
